# Applies the "Long Stayers (Local Settings)" guidance-text edit to the
# EvaChecks workbook:
#   - On "DQ Checks": rename the existing "Long Stayers" issue to
#     "Long Stayers (Local Settings)", drop its old Notes text, and add a
#     brand-new "Long Stayers" row underneath with fresh guidance text.
#   - Mirror the same change on "All Checks", where the new row is
#     inserted in the middle of the table (so everything below it shifts
#     down by one row) and the AutoFilter / _FilterDatabase range grows
#     to match.

$wb = $excel.ActiveWorkbook

$guidanceText = 'This household has been housed in your project for a relatively long time compared to enrollments into the same project type in the rest of your system. If they have exited, please enter an Exit Date, otherwise consider using Move On Assistance funds. If they need to remain in the project, leave everything as is.'

# ---------------------------------------------------------------------
# 1. "DQ Checks" sheet
# ---------------------------------------------------------------------
$wsDQ = $wb.Worksheets.Item("DQ Checks")

# Row 73: "Long Stayers" -> "Long Stayers (Local Settings)"; old Notes
# ("Includes Services Only, Other, and Day Shelter") is removed.
$wsDQ.Range("C73").Value = "Long Stayers (Local Settings)"
$wsDQ.Range("F73").ClearContents()
$wsDQ.Range("F73").Style = "Normal"

# New row 75 (appended after the current last row, 74) with the new
# "Long Stayers" guidance issue.
$wsDQ.Range("A75").Value = "dq"
$wsDQ.Range("B75").Value = "Warning"
$wsDQ.Range("C75").Value = "Long Stayers"
$wsDQ.Range("D75").Value = "3.10 - Project Start Date"
$wsDQ.Range("E75").Value = $guidanceText
$wsDQ.Range("G75").Value = 104

# ---------------------------------------------------------------------
# 2. "All Checks" sheet
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Checks")

# Row 98: same text change as DQ Checks row 73.
$wsAll.Range("C98").Value = "Long Stayers (Local Settings)"
$wsAll.Range("F98").ClearContents()
$wsAll.Range("F98").Style = "Normal"

# Insert a new row at 100 (shifts the existing rows 100-104 down to
# 101-105, carrying their formatting/row-height with them) and fill it
# with the new "Long Stayers" guidance issue.
$wsAll.Rows.Item(100).Insert()

$wsAll.Range("A100").Value = "dq"
$wsAll.Range("B100").Value = "Warning"
$wsAll.Range("C100").Value = "Long Stayers"
$wsAll.Range("D100").Value = "3.10 - Project Start Date"
$wsAll.Range("E100").Value = $guidanceText
$wsAll.Range("F100").ClearContents()
$wsAll.Range("F100").Style = "Normal"
$wsAll.Range("G100").Value = 104

# Re-apply the AutoFilter so its range grows from A1:G104 to A1:G105.
$wsAll.AutoFilterMode = $false
$wsAll.Range("A1:G105").AutoFilter()

# Update the hidden _FilterDatabase defined name to match.
$fdb = $wb.Names.Item("All Checks!_FilterDatabase")
$fdb.RefersTo = "='All Checks'!`$A`$1:`$G`$105"

# ---------------------------------------------------------------------
# 3. Selection / view bookkeeping (cosmetic, best effort)
# ---------------------------------------------------------------------
$wsDQ.Activate()
$wsDQ.Range("C84").Select()

$wsAll.Activate()
$wsAll.Range("F97").Select()
